# Add 4 new rows (173-176) of betting data to the "bets" sheet, matching the
# author's upload of new TENIS DE MESA / LIGA PRO bets dated 2023-12-04/05.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

function Add-BetRow($ws, $Row, $ColA, $ColB, $ColC, $ColD, $ColE, $ColF, $ColG, $ColH, $ColI, $ColJ, $ColK, $ColL) {
    $prev = $Row - 1

    $ws.Cells.Item($Row, 1).Value = $ColA
    $ws.Cells.Item($Row, 2).Value = $ColB

    $ws.Cells.Item($Row, 3).Value = $ColC
    [void]$ws.Cells.Item($prev, 3).Copy()
    [void]$ws.Cells.Item($Row, 3).PasteSpecial(-4122)

    $ws.Cells.Item($Row, 4).Value = $ColD
    [void]$ws.Cells.Item($prev, 4).Copy()
    [void]$ws.Cells.Item($Row, 4).PasteSpecial(-4122)

    $ws.Cells.Item($Row, 5).Value = $ColE
    $ws.Cells.Item($Row, 6).Value = $ColF

    $ws.Cells.Item($Row, 7).Formula = $ColG
    $ws.Cells.Item($Row, 8).Value = $ColH
    $ws.Cells.Item($Row, 9).Formula = $ColI

    $ws.Cells.Item($Row, 10).Value = $ColJ
    $ws.Cells.Item($Row, 11).Value = $ColK

    $ws.Cells.Item($Row, 12).Formula = $ColL
    [void]$ws.Cells.Item($prev, 12).Copy()
    [void]$ws.Cells.Item($Row, 12).PasteSpecial(-4122)
}

Add-BetRow $ws 173 172 155 45264 "2023-12-04" 1 1.085 "=I172" 88 "=G173+H173" "TENIS DE MESA" "MASTERS" '=ROUND((I173/$G$31-1)*100, 3)+$L$29'

Add-BetRow $ws 174 173 156 45264 "2023-12-04" 1 1.0640000000000001 "=I173" 51 "=G174+H174" "TENIS DE MESA" "MASTERS" '=ROUND((I174/$G$31-1)*100, 3)+$L$29'

Add-BetRow $ws 175 174 157 45265 "2023-12-05" 1 1.32 "=I174" 272 "=G175+H175" "TENIS DE MESA" "MASTERS" '=ROUND((I175/$G$31-1)*100, 3)+$L$29'

Add-BetRow $ws 176 175 158 45265 "2023-12-05" 1 1.1100000000000001 "=I175" 79 "=G176+H176" "TENIS DE MESA" "LIGA PRO" '=ROUND((I176/$G$31-1)*100, 3)+$L$29'

# Restore the selection / scroll position recorded in the saved view.
[void]$ws.Range("K174").Select()

Write-Host "Added rows 173-176 to bets sheet"
